$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$emdash = [char]0x2014
$rsquo = [char]0x2019

$ws.Range("A2").Value = "('Frost Titan', ['{4}{U}{U}', 'Creature $emdash Giant', 'Whenever Frost Titan becomes the target of a spell or ability an opponent controls, counter that spell or ability unless its controller pays {2}.', 'Whenever Frost Titan enters the battlefield or attacks, tap target permanent. It doesn${rsquo}t untap during its controller${rsquo}s next untap step.', '6/6'])"

$ws.Range("A3").Value = "('Grave Titan', ['{4}{B}{B}', 'Creature $emdash Giant', 'Deathtouch', 'Whenever Grave Titan enters the battlefield or attacks, create two 2/2 black Zombie creature tokens.', '6/6'])"

$ws.Range("A4").Value = "('Inferno Titan', ['{4}{R}{R}', 'Creature $emdash Giant', '{R}: Inferno Titan gets +1/+0 until end of turn.', 'Whenever Inferno Titan enters the battlefield or attacks, it deals 3 damage divided as you choose among one, two, or three targets.', '6/6'])"

$ws.Range("A5:A19").Clear()
